$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new "Citation" column after "Article" (old B..H shift to C..I)
$ws.Columns("B:B").Insert()

# 2) Insert a blank spacer column before "n_studies" (old F..H, now at G..I,
#    shift to H..J), matching the new layout's blank column G
$ws.Columns("G:G").Insert()
$ws.Range("G3").Clear()

# 3) Header for the new column
$ws.Range("B1").Value = "Citation"

# 4) Row 2 (RPP): move the long citation that used to live in A2 into the new
#    B2 "Citation" cell, and replace A2 with the short study-name label "RPP"
$rppCitation = $ws.Range("A2").Value2
$ws.Range("B2").Value = $rppCitation
$ws.Range("A2").Value = "RPP"

# 5) Row 3 (Many labs 1): add the Klein et al. citation
$ws.Range("B3").Value = "Klein, R. A., Ratliff, K. A., Vianello, M., Adams, R. B., Bahník, Š., Bernstein, M. J., . . . Nosek, B. A. (2014). Investigating Variation in Replicability. Social Psychology, 45(3), 142-152. doi:10.1027/1864-9335/a000178"

# 6) Row 5 (Many Labs 3): add the Ebersole et al. citation
$ws.Range("B5").Value = "Ebersole, C. R., Atherton, O. E., Belanger, A. L., Skulborstad, H. M., Allen, J. M., Banks, J. B., . . . Nosek, B. A. (2016). Many Labs 3: Evaluating participant pool quality across the academic semester via replication. Journal of Experimental Social Psychology, 67, 68-82. doi:https://doi.org/10.1016/j.jesp.2015.10.012"

# 7) Re-create the hyperlinks that the column insert left stranded on the old
#    (now-empty) column letters, pointing them at the cells that now hold the
#    linked text. For A2/B2 seed the cached display text with the URL first
#    (matching the original file's hyperlink metadata) then restore the real
#    citation text afterwards.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "http://www.sciencemag.org/cgi/content/full/349/6251/aac4716?ijkey=1xgFoCnpLswpk&keytype=ref&siteid=sci", [System.Type]::Missing, [System.Type]::Missing, "http://www.sciencemag.org/cgi/content/full/349/6251/aac4716?ijkey=1xgFoCnpLswpk&keytype=ref&siteid=sci") | Out-Null
$ws.Range("A2").Value = "RPP"
$ws.Hyperlinks.Add($ws.Range("B2"), "http://www.sciencemag.org/cgi/content/full/349/6251/aac4716?ijkey=1xgFoCnpLswpk&keytype=ref&siteid=sci", [System.Type]::Missing, [System.Type]::Missing, "http://www.sciencemag.org/cgi/content/full/349/6251/aac4716?ijkey=1xgFoCnpLswpk&keytype=ref&siteid=sci") | Out-Null
$ws.Range("B2").Value = $rppCitation
$ws.Hyperlinks.Add($ws.Range("C4"), "https://osf.io/8cd4r/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://osf.io/ct89g/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://osf.io/pfdyw/") | Out-Null

# 8) Add the n_studies total row
$ws.Range("H10").Formula = "=SUM(H2:H9)"

# 9) Match the saved selection in the source workbook
$ws.Range("I10").Select()

Write-Host "edit complete"
